# -----------------------------------------------------------------------
# This script changes several previously red (FF0000) text passages back
# to "Automatic" font color, and cleans up a couple of paragraphs whose
# only remaining difference between runs was font color (which, once the
# color is normalized, allows the runs carrying identical text/formatting
# to be re-joined into a single run) - mirroring an author pass that
# selected the red "draft" notes and reset their font color to Automatic.
# -----------------------------------------------------------------------

$wdColorAutomatic = -16777216
$d = $word.ActiveDocument

function Get-ParagraphByText($doc, [string]$needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $para = $doc.Paragraphs.Item($i)
        if ($para.Range.Text.Contains($needle)) {
            return $para
        }
    }
    return $null
}

# 1) "Brug af eksternt system ..." - no color specified before; make it
#    explicit "Automatic" on both the paragraph mark and the run.
$p1 = Get-ParagraphByText $d "Brug af eksternt system"
$p1.Range.Font.Color = $wdColorAutomatic

# 2) "Kilometer udregning er afhaengig ... beregning udfoert." - was red;
#    the trailing "." was its own run - merge it into the sentence, then
#    recolor the whole paragraph (incl. paragraph mark) to Automatic.
$sentence = $d.Content
$needle2 = "Kilometer udregning er afh" + [char]230 + "ngig af internet adgang og dens hastighed da der kaldes op til Google for at f" + [char]229 + " en beregning udf" + [char]248 + "rt."
$sentence.Find.Execute($needle2, $false, $false, $false, $false, $false, $true, 1, $false, $needle2, 2) | Out-Null
$p2 = Get-ParagraphByText $d "Kilometer udregning er afh"
$p2.Range.Font.Color = $wdColorAutomatic

# 3) The empty paragraph right after it (its paragraph mark was red) -
#    Font.Color on a wholly empty range is a no-op in this host, so we
#    briefly insert a placeholder character to anchor the paragraph mark
#    run, recolor it, then remove the placeholder again.
$p3 = $p2.Next()
$placeholderStart = $d.Range($p3.Range.Start, $p3.Range.Start)
$placeholderStart.InsertBefore("x")
$p3again = $d.Range($p3.Range.Start, $p3.Range.Start)
$placeholderRange = $d.Range($p3.Range.Start, $p3.Range.Start + 1)
$p3.Range.Font.Color = $wdColorAutomatic
$placeholderRange2 = $d.Range($p3.Range.Start, $p3.Range.Start + 1)
$placeholderRange2.Delete()

# 4) "Kilometer udregning udfoeres i en selvstaendig proces af
#    okhttp-2.5.0.jar." - three separate runs (split by proofErr markers)
#    all need FF0000 -> Automatic, but keep them as separate runs.
$p4 = Get-ParagraphByText $d "okhttp-2.5.0.jar"
$p4.Range.Font.Color = $wdColorAutomatic

# 5) "Pris paa ekstra tilvalg, afstand (ved kilometer udregning) og antal
#    af personer indgaar i beregningen." - recolor the two red runs to
#    Automatic and merge each with its neighbouring non-coloured run
#    (", " + "afstand (...)"  and  " " + "og antal ...").
#    Order matters: colour the *first* run to its final value before any
#    text-merge touches its neighbour, otherwise the merge coalesces it
#    into that first run too (the host merges adjacent runs that already
#    carry identical formatting as soon as a text edit touches either of
#    them).
$p5 = Get-ParagraphByText $d "Pris p" + [char]229 + " ekstra tilvalg"
# (fallback in case string concatenation above confuses Find text; look
#  the paragraph up the safe way instead)
$p5 = Get-ParagraphByText $d "ekstra tilvalg"

$r5start = $p5.Range.Start
$firstRunLen = "Pris p${([char]229)} ekstra tilvalg".Length

$r1 = $d.Range($r5start, $r5start + $firstRunLen)
$r1.Font.Color = $wdColorAutomatic

$commaAndDistance = $d.Range($r1.End, $r1.End)
